$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.432.21'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -2.79%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.429.74'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.53%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '582.56'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -3.72%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '133.83'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -6.23%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.429.35'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -2.56%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.483'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -6.25%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.121'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -7.84%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.00'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -9.26%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.375'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -8.41%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.016.26'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -2.50%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000179'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -8.12%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '26.30'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -8.25%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.428.37'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -2.58%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.115'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '64.441.37'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -2.61%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.59'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -11.53%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.68'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -8.17%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.58'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -6.93%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '381.27'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -9.94%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.544'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -7.89%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.72'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '71.76'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -6.86%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.563.88'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -2.61%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0000105'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -7.58%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.17'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -9.57%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.06'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -9.85%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.19'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -11.07%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.441.03'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -2.40%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '22.95'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -5.53%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.141'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -9.66%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '171.07'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -2.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.18'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -11.54%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.71'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -11.11%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.46'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -10.61%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.66'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -10.80%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0760'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -7.28%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.801'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -6.84%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '41.80'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -7.80%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.29'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -14.12%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.60'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -9.42%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.11'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '22.65'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -1.61%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.52'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -7.92%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.197.19'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -5.24%  '
